$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pros1"
$ws.Range("C2").Value = "Tyro3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 18.741769
$ws.Range("H2").Value = 56.225307
$ws.Range("I2").Value = 0.2218531826860132
$ws.Range("J2").Value = 0.2218531826860132
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.448064
$ws.Range("N2").Value = 1.344192
$ws.Range("O2").Value = 0.113372348992564
$ws.Range("P2").Value = 0.113372348992564
$ws.Range("Q2").Value = 8.397511985216001
$ws.Range("R2").Value = 75.577607866944
$ws.Range("S2").Value = 0.02515201645258976
$ws.Range("T2").Value = 0.02515201645258976

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pros1"
$ws.Range("C3").Value = "Tyro3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 18.741769
$ws.Range("H3").Value = 56.225307
$ws.Range("I3").Value = 0.2218531826860132
$ws.Range("J3").Value = 0.2218531826860132
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.951476
$ws.Range("N3").Value = 8.854428
$ws.Range("O3").Value = 0.7468035082380574
$ws.Range("P3").Value = 0.7468035082380574
$ws.Range("Q3").Value = 55.315881401044
$ws.Range("R3").Value = 497.842932609396
$ws.Range("S3").Value = 0.1656807351436933
$ws.Range("T3").Value = 0.1656807351436933

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pros1"
$ws.Range("C4").Value = "Tyro3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 18.741769
$ws.Range("H4").Value = 56.225307
$ws.Range("I4").Value = 0.2218531826860132
$ws.Range("J4").Value = 0.2218531826860132
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5526053333333333
$ws.Range("N4").Value = 1.657816
$ws.Range("O4").Value = 0.1398241427693786
$ws.Range("P4").Value = 0.1398241427693786
$ws.Range("Q4").Value = 10.35680150550133
$ws.Range("R4").Value = 93.211213549512
$ws.Range("S4").Value = 0.03102043108973014
$ws.Range("T4").Value = 0.03102043108973014

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pros1"
$ws.Range("C5").Value = "Tyro3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 53.77230066666667
$ws.Range("H5").Value = 161.316902
$ws.Range("I5").Value = 0.6365224138259964
$ws.Range("J5").Value = 0.6365224138259964
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.448064
$ws.Range("N5").Value = 1.344192
$ws.Range("O5").Value = 0.113372348992564
$ws.Range("P5").Value = 0.113372348992564
$ws.Range("Q5").Value = 24.09343212590933
$ws.Range("R5").Value = 216.840889133184
$ws.Range("S5").Value = 0.07216404124187013
$ws.Range("T5").Value = 0.07216404124187013

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pros1"
$ws.Range("C6").Value = "Tyro3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 53.77230066666667
$ws.Range("H6").Value = 161.316902
$ws.Range("I6").Value = 0.6365224138259964
$ws.Range("J6").Value = 0.6365224138259964
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.951476
$ws.Range("N6").Value = 8.854428
$ws.Range("O6").Value = 0.7468035082380574
$ws.Range("P6").Value = 0.7468035082380574
$ws.Range("Q6").Value = 158.7076548824507
$ws.Range("R6").Value = 1428.368893942056
$ws.Range("S6").Value = 0.4753571717174107
$ws.Range("T6").Value = 0.4753571717174107

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pros1"
$ws.Range("C7").Value = "Tyro3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 53.77230066666667
$ws.Range("H7").Value = 161.316902
$ws.Range("I7").Value = 0.6365224138259964
$ws.Range("J7").Value = 0.6365224138259964
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5526053333333333
$ws.Range("N7").Value = 1.657816
$ws.Range("O7").Value = 0.1398241427693786
$ws.Range("P7").Value = 0.1398241427693786
$ws.Range("Q7").Value = 29.71486013400355
$ws.Range("R7").Value = 267.433741206032
$ws.Range("S7").Value = 0.08900120086671559
$ws.Range("T7").Value = 0.08900120086671559

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pros1"
$ws.Range("C8").Value = "Tyro3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 11.964182
$ws.Range("H8").Value = 35.892546
$ws.Range("I8").Value = 0.1416244034879904
$ws.Range("J8").Value = 0.1416244034879904
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.448064
$ws.Range("N8").Value = 1.344192
$ws.Range("O8").Value = 0.113372348992564
$ws.Range("P8").Value = 0.113372348992564
$ws.Range("Q8").Value = 5.360719243648
$ws.Range("R8").Value = 48.246473192832
$ws.Range("S8").Value = 0.01605629129810416
$ws.Range("T8").Value = 0.01605629129810415

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pros1"
$ws.Range("C9").Value = "Tyro3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 11.964182
$ws.Range("H9").Value = 35.892546
$ws.Range("I9").Value = 0.1416244034879904
$ws.Range("J9").Value = 0.1416244034879904
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.951476
$ws.Range("N9").Value = 8.854428
$ws.Range("O9").Value = 0.7468035082380574
$ws.Range("P9").Value = 0.7468035082380574
$ws.Range("Q9").Value = 35.311996032632
$ws.Range("R9").Value = 317.807964293688
$ws.Range("S9").Value = 0.1057656013769534
$ws.Range("T9").Value = 0.1057656013769534

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pros1"
$ws.Range("C10").Value = "Tyro3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 11.964182
$ws.Range("H10").Value = 35.892546
$ws.Range("I10").Value = 0.1416244034879904
$ws.Range("J10").Value = 0.1416244034879904
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.5526053333333333
$ws.Range("N10").Value = 1.657816
$ws.Range("O10").Value = 0.1398241427693786
$ws.Range("P10").Value = 0.1398241427693786
$ws.Range("Q10").Value = 6.611470782170666
$ws.Range("R10").Value = 59.50323703953599
$ws.Range("S10").Value = 0.01980251081293285
$ws.Range("T10").Value = 0.01980251081293285
